# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
#   on every sheet that reports it (Overview: E2/F2, zh-cn: C2, de-de: C2)
# - Narrow the now-shorter "Status" columns to match the re-generated
#   report's auto-fit width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# A ColumnWidth of 12.5 is the value this host's width quantization maps
# to the narrowed, auto-fit column width used by the regenerated report.
$newColumnWidth = 12.5

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Columns.Item(5).ColumnWidth = $newColumnWidth
$ws1.Columns.Item(6).ColumnWidth = $newColumnWidth

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = $newStatus
$ws2.Columns.Item(3).ColumnWidth = $newColumnWidth

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = $newStatus
$ws3.Columns.Item(3).ColumnWidth = $newColumnWidth
